$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 01:52"

# --- Update Estados Unidos (row 4) ---
$ws.Range("B4").Value = 395612
$ws.Range("C4").Value = 28608
$ws.Range("E4").Value = 361148
$ws.Range("G4").Value = 1919
$ws.Range("H4").Value = 12790

# --- Update Australia (row 25) ---
$ws.Range("B25").Value = 5988
$ws.Range("C25").Value = 93
$ws.Range("E25").Value = 3392
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 49

# --- Ecuador overtakes Malasia/Japon/Filipinas in the ranking (rows 34-37) ---
# Row 34 keeps its row position but now shows Ecuador with updated figures
$ws.Range("A34").Value = "Ecuador"
$ws.Range("B34").Value = 3995
$ws.Range("C34").Value = 248
$ws.Range("D34").Value = 140
$ws.Range("E34").Value = 3635
$ws.Range("F34").Value = 156
$ws.Range("G34").Value = 29
$ws.Range("H34").Value = 220

# Row 35 now shows Malasia (its figures are unchanged, only its rank shifted)
$ws.Range("A35").Value = "Malasia"
$ws.Range("B35").Value = 3963
$ws.Range("C35").Value = 170
$ws.Range("D35").Value = 1321
$ws.Range("E35").Value = 2579
$ws.Range("F35").Value = 92
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 63

# Row 36 now shows Japon (figures unchanged)
$ws.Range("A36").Value = "Japon"
$ws.Range("B36").Value = 3906
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 592
$ws.Range("E36").Value = 3222
$ws.Range("F36").Value = 79
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 92

# Row 37 now shows Filipinas (figures unchanged)
$ws.Range("A37").Value = "Filipinas"
$ws.Range("B37").Value = 3764
$ws.Range("C37").Value = 104
$ws.Range("D37").Value = 84
$ws.Range("E37").Value = 3503
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 14
$ws.Range("H37").Value = 177

# --- Update Argentina (row 53) ---
$ws.Range("B53").Value = 1715
$ws.Range("C53").Value = 87
$ws.Range("E53").Value = 1317
$ws.Range("G53").Value = 7
$ws.Range("H53").Value = 60

# --- Update Camerun (row 77) ---
$ws.Range("B77").Value = 685
$ws.Range("C77").Value = 27
$ws.Range("D77").Value = 60
$ws.Range("E77").Value = 616

# --- Guayana Francesa overtakes Aruba in the ranking (rows 135-136) ---
# Row 135 keeps its row position but now shows Guayana Francesa with updated figures
$ws.Range("A135").Value = "Guayana Francesa"
$ws.Range("B135").Value = 77
$ws.Range("C135").Value = 5
$ws.Range("D135").Value = 34
$ws.Range("E135").Value = 43
$ws.Range("F135").Value = 1

# Row 136 now shows Aruba (figures unchanged, only its rank shifted)
$ws.Range("A136").Value = "Aruba"
$ws.Range("B136").Value = 74
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 14
$ws.Range("E136").Value = 60
$ws.Range("F136").Value = 0
